# Arabic translation pass for "Email 2-1 [TEMPLATE] Partner email – reminder to RSVP.docx"
#
# Each replacement is scoped to the specific paragraph Range it belongs to
# (rather than $d.Content) so that Find/Replace only ever touches the run(s)
# inside that paragraph. This prevents unrelated matches elsewhere in the
# document (e.g. duplicate "English"/"live chat"/", " strings) from being
# altered, and keeps runs with different formatting (highlights, hyperlinks,
# etc.) intact instead of being merged into a single run.

$d = $word.ActiveDocument

function Replace-InParagraph([int]$index, [string]$find, [string]$replace) {
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Paragraph 1: "English / Portuguese / French / Thai / Vietnamese / Spanish"
Replace-InParagraph 1 "English" "الإنجليزية"
Replace-InParagraph 1 " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

# Paragraph 3: "English" (heading style)
Replace-InParagraph 3 "English" "الإنجليزية"

# Paragraph 5: "Brief:" -> only the "Brief" word changes, ":" run stays
Replace-InParagraph 5 "Brief" "المضمون"

# Paragraph 6: brief description sentence, only the tail changes
Replace-InParagraph 6 "An email to partners in the target country who haven’t RSVPed to remind them to send the RSVP. It will be sent via customer.io" "An email to partners in the target country who haven’t RSVPed to remind them to send the RSVP. سيتم إرسالها عبر customer.io"

# Paragraph 8: "Target audience:" -> only "Target audience" changes
Replace-InParagraph 8 "Target audience" "الجمهور المستهدف"

# Paragraph 15: "Don’t delay! Book your spot today!"
Replace-InParagraph 15 "Don’t delay! Book your spot today!" "لا تتأخر! احجز مكانك اليوم!"

# Paragraph 17: "Hi [PARTNER NAME], " -> "Hi " and ", " runs change; highlighted
# "[PARTNER NAME]" run stays as-is.
Replace-InParagraph 17 "Hi " "مرحبًا  "
Replace-InParagraph 17 ", " ",، "

# Paragraph 25: "If you have any questions, please contact us via live chat or WhatsApp."
Replace-InParagraph 25 "If you have any questions, please contact us via " "إذا كانت لديك أي أسئلة، فاتصل بنا:  "
Replace-InParagraph 25 "live chat" "الدردشة الحية"

# Paragraph 26: "If you have any questions, please contact your country manager, [NAME], ..."
Replace-InParagraph 26 "If you have any questions, please contact your country manager, " "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمدير بلدك  "

# Paragraph 27: "We look forward to seeing you at [EVENT NAME]! "
Replace-InParagraph 27 "We look forward to seeing you at [EVENT NAME]! " "نتطلع إلى رؤيتك في [EVENT NAME]! "

# Paragraph 39: "If you have any questions, please contact your country manager:"
Replace-InParagraph 39 "If you have any questions, please contact your country manager:" "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمديرك الإقليمي:"

Write-Host "Arabic translation replacements applied."
